# "incluida simulacao do compras" - add simulated purchase order item rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (row 1 header): A=seq, B=seq_pedido, C=seq_item, D=descricao_item,
#                          E=valor_unitario, F=quantidade, G=desconto
$rows = @(
    @(1, 1, 10, "Item 10", 2,    10, 0),
    @(2, 1, 11, "Item 11", 45,   2,  5),
    @(3, 1, 12, "Item 12", 8.5,  5,  0)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $r = $r + 1
}

$ws.Range("A1").Select()
